$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 20) with a new forecast vector entry, matching the
# pattern of the existing rows (date, year, yoy value, next year, next yoy).
$ws.Range("A20").Value = 45986
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = -2.06674933094535
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = -0.9969640812590996
